$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("C25").Value = 1036.8
$ws1.Range("C38").Value = 1036.8
$ws1.Range("I38").Value = 273.6
$ws1.Range("I44").Value = 369
$ws1.Range("C57").Value = "2 de 55"
$ws1.Range("I57").Value = "3 de 55"

# --- Sheet "VENTA MENSUAL" ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F25").Value = 4515.21
$ws2.Range("F38").Value = 2718.72
$ws2.Range("F44").Value = 858.12
$ws2.Range("F57").Value = 40960.37

# --- Sheet "CUMPLIMIENTO MENSUAL" ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D2").Value = 1974.07
$ws3.Range("E2").Value = 7996.27304517915
$ws3.Range("F2").Value = 0.1979941904761743

$ws3.Range("D8").Value = 730.5599999999999
$ws3.Range("E8").Value = 269.4400000000001
$ws3.Range("F8").Value = 0.73056

$ws3.Range("D19").Value = 40960.37
$ws3.Range("E19").Value = 76479.32064517915
$ws3.Range("F19").Value = 0.3487779112408741

# Column F width 25 -> 24 (stored XML width). ColumnWidth property uses a
# different unit than the stored "width" attribute (offset of 5/6 for this
# workbook's default font), so 23.166666666666668 yields a stored width of 24.
$ws3.Columns.Item(6).ColumnWidth = 23.166666666666668
